$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 233.46153
$ws.Range("I11").Value = 233.46153
$ws.Range("K11").Value = 233.46153
$ws.Range("M11").Value = -93.46153000000001
$ws.Range("H18").Value = 899.6667
$ws.Range("I18").Value = 899.6667
$ws.Range("K18").Value = 899.6667
$ws.Range("M18").Value = -615.6667
$ws.Range("H74").Value = 59177
$ws.Range("I74").Value = 4074.125
$ws.Range("K74").Value = 4074.125
$ws.Range("M74").Value = -3138.125
$ws.Range("H77").Value = 59177
$ws.Range("I77").Value = 4074.125
$ws.Range("K77").Value = 20370.625
$ws.Range("M77").Value = -15690.625
$ws.Range("H93").Value = 43466.5
$ws.Range("J93").Value = 43466.5
$ws.Range("L93").Value = 43466.5
$ws.Range("N93").Value = -48458.5
$ws.Range("H132").Value = 1588.2759
$ws.Range("I132").Value = 1588.2759
$ws.Range("K132").Value = 4764.8277
$ws.Range("M132").Value = -2234.8277
$ws.Range("H137").Value = 3849.1667
$ws.Range("I137").Value = 1200
$ws.Range("K137").Value = 3600
$ws.Range("M137").Value = -1050
$ws.Range("H141").Value = 4475.9
$ws.Range("I141").Value = 4475.9
$ws.Range("K141").Value = 13427.7
$ws.Range("M141").Value = -8247.699999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1806.4615
$ws.Range("I45").Value = 1748.6666
$ws.Range("K45").Value = 1748.6666
$ws.Range("M45").Value = -1371.6666
$ws.Range("H61").Value = 2403.25
$ws.Range("I61").Value = 2403.25
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2403.25
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2191.25
$ws.Range("N61").ClearContents()
$ws.Range("H110").Value = 1888.7142
$ws.Range("I110").Value = 1741.091
$ws.Range("K110").Value = 1741.091
$ws.Range("M110").Value = 303.9090000000001
$ws.Range("H132").Value = 1977.6428
$ws.Range("I132").Value = 1975.9231
$ws.Range("K132").Value = 5927.7693
$ws.Range("M132").Value = -3397.7693
$ws.Range("H136").Value = 2403.25
$ws.Range("I136").Value = 2403.25
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7209.75
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4659.75
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 983
$ws.Range("I107").Value = 962.125
$ws.Range("K107").Value = 962.125
$ws.Range("M107").Value = 957.875
$ws.Range("H134").Value = 3029.7917
$ws.Range("I134").Value = 3029.7917
$ws.Range("K134").Value = 9089.375100000001
$ws.Range("M134").Value = -6554.375100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2788.9443
$ws.Range("I31").Value = 2291.4348
$ws.Range("J31").Value = 3669.1538
$ws.Range("K31").Value = 2291.4348
$ws.Range("L31").Value = 3669.1538
$ws.Range("M31").Value = -1996.4348
$ws.Range("N31").Value = -4259.1538
$ws.Range("H34").Value = 2788.9443
$ws.Range("I34").Value = 2291.4348
$ws.Range("J34").Value = 3669.1538
$ws.Range("K34").Value = 2291.4348
$ws.Range("L34").Value = 3669.1538
$ws.Range("M34").Value = -2089.4348
$ws.Range("N34").Value = -4073.1538
$ws.Range("H58").Value = 1721.5454
$ws.Range("I58").Value = 1677.5555
$ws.Range("K58").Value = 1677.5555
$ws.Range("M58").Value = -1474.5555
$ws.Range("H136").Value = 1721.5454
$ws.Range("I136").Value = 1677.5555
$ws.Range("K136").Value = 5032.666499999999
$ws.Range("M136").Value = -2482.666499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2250
$ws.Range("I64").Value = 1500
$ws.Range("J64").Value = 2500
$ws.Range("K64").Value = 4500
$ws.Range("L64").Value = 7500
$ws.Range("M64").Value = -4230
$ws.Range("N64").Value = -8040
$ws.Range("H67").Value = 2250
$ws.Range("I67").Value = 1500
$ws.Range("J67").Value = 2500
$ws.Range("K67").Value = 4500
$ws.Range("L67").Value = 7500
$ws.Range("M67").Value = -3564
$ws.Range("N67").Value = -9372
$ws.Range("H80").Value = 7070
$ws.Range("I80").Value = 140
$ws.Range("J80").Value = 14000
$ws.Range("K80").Value = 420
$ws.Range("L80").Value = 42000
$ws.Range("M80").Value = 516
$ws.Range("N80").Value = -43872
$ws.Range("H83").Value = 7070
$ws.Range("I83").Value = 140
$ws.Range("J83").Value = 14000
$ws.Range("K83").Value = 1260
$ws.Range("L83").Value = 126000
$ws.Range("M83").Value = 3420
$ws.Range("N83").Value = -135360
$ws.Range("H99").Value = 2391.6667
$ws.Range("I99").Value = 2391.6667
$ws.Range("K99").Value = 7175.000100000001
$ws.Range("M99").Value = -4929.000100000001
$ws.Range("H113").Value = 1378.3
$ws.Range("I113").Value = 997.6667
$ws.Range("J113").Value = 1949.25
$ws.Range("K113").Value = 2993.0001
$ws.Range("L113").Value = 5847.75
$ws.Range("M113").Value = -823.0001000000002
$ws.Range("N113").Value = -10187.75
$ws.Range("H138").Value = 5848.6
$ws.Range("J138").Value = 6853
$ws.Range("L138").Value = 20559
$ws.Range("N138").Value = -30839

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3435.25
$ws.Range("I80").Value = 1500
$ws.Range("J80").Value = 3711.7144
$ws.Range("K80").Value = 1500
$ws.Range("L80").Value = 3711.7144
$ws.Range("M80").Value = -502
$ws.Range("N80").Value = -5707.7144
$ws.Range("H83").Value = 3435.25
$ws.Range("I83").Value = 1500
$ws.Range("J83").Value = 3711.7144
$ws.Range("K83").Value = 7500
$ws.Range("L83").Value = 18558.572
$ws.Range("M83").Value = -2508
$ws.Range("N83").Value = -28542.572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1679.3
$ws.Range("I22").Value = 1310.3334
$ws.Range("K22").Value = 1310.3334
$ws.Range("M22").Value = -1015.3334
$ws.Range("H27").Value = 1679.3
$ws.Range("I27").Value = 1310.3334
$ws.Range("K27").Value = 1310.3334
$ws.Range("M27").Value = -1203.3334
$ws.Range("H122").Value = 2951.5
$ws.Range("I122").Value = 2930.7144
$ws.Range("K122").Value = 8792.143199999999
$ws.Range("M122").Value = -6342.143199999999
$ws.Range("H132").Value = 4812.5
$ws.Range("I132").Value = 2392.8572
$ws.Range("J132").Value = 8200
$ws.Range("K132").Value = 7178.571599999999
$ws.Range("L132").Value = 24600
$ws.Range("M132").Value = -4648.571599999999
$ws.Range("N132").Value = -29660

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H122").Value = 2079.4
$ws.Range("I122").Value = 1849.25
$ws.Range("K122").Value = 5547.75
$ws.Range("M122").Value = -3097.75
$ws.Range("H132").Value = 2189.5833
$ws.Range("I132").Value = 2189.5833
$ws.Range("K132").Value = 6568.749899999999
$ws.Range("M132").Value = -4038.749899999999
$ws.Range("H136").Value = 1750.7222
$ws.Range("I136").Value = 1407.4615
$ws.Range("J136").Value = 2643.2
$ws.Range("K136").Value = 4222.3845
$ws.Range("L136").Value = 7929.599999999999
$ws.Range("M136").Value = -1672.3845
$ws.Range("N136").Value = -13029.6
